$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update column F values for specific rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1211
$ws1.Range("F4").Value = 16890
$ws1.Range("F5").Value = 35
$ws1.Range("F6").Value = 1647
$ws1.Range("F7").Value = 69
$ws1.Range("F9").Value = 391
$ws1.Range("F12").Value = 11708
$ws1.Range("F13").Value = 28
$ws1.Range("F14").Value = 1375
$ws1.Range("F15").Value = 4649
$ws1.Range("F16").Value = 464
$ws1.Range("F17").Value = 4
$ws1.Range("F20").Value = 897
$ws1.Range("F23").Value = 21
$ws1.Range("F24").Value = 5213

# Sheet "全部类型" (sheet4) - same underlying events, rows offset by +1
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1211
$ws4.Range("F5").Value = 16890
$ws4.Range("F6").Value = 35
$ws4.Range("F7").Value = 1647
$ws4.Range("F8").Value = 69
$ws4.Range("F10").Value = 391
$ws4.Range("F15").Value = 11708
$ws4.Range("F16").Value = 28
$ws4.Range("F17").Value = 1375
$ws4.Range("F18").Value = 4649
$ws4.Range("F19").Value = 464
$ws4.Range("F20").Value = 4
$ws4.Range("F23").Value = 897
$ws4.Range("F26").Value = 21
$ws4.Range("F27").Value = 5213
